$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 3
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -5
$ws.Range("F12").Value = 0
